$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Four new labelled "summary" rows (14-17) below the per-instance data table:
#   Average / Worst of the SW(S*)/SW(OPT) ratio (col N) and the
#   SC(S*)/SC(OPT) ratio (col Z).
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style the 4 new "B" values (bold, 12pt, vertically centred). The combined
# font/alignment is built once on a scratch cell and copied across so the
# whole range picks up a single new style entry instead of one per property
# assignment.
$scratch = $ws.Range("Z30")
$scratch.Font.Bold = $true
$scratch.Font.Size = 12
$scratch.VerticalAlignment = -4108   # xlVAlignCenter
$scratch.Copy()
$ws.Range("B14:B17").PasteSpecial(-4122)   # xlPasteFormats
$scratch.Clear()

# ---------------------------------------------------------------------------
# New row 12: average of the k column (J) directly under the data table.
# ---------------------------------------------------------------------------
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# Page setup tweak that came along with this edit.
$ws.PageSetup.PaperSize = 9          # xlPaperA4
$ws.PageSetup.Orientation = 1        # xlPortrait

# Leave the selection where the author left it when saving.
$ws.Range("J12").Select() | Out-Null
